$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing rows shift down by one
$ws.Range("A1").EntireRow.Insert()

# New header row
$ws.Range("A1").Value = "Raum"
$ws.Range("B1").Value = "Kapazität"

# New capacity values in column B for the (now shifted) data rows
$ws.Range("B2").Value = 15
$ws.Range("B3").Value = 20
$ws.Range("B4").Value = 20
$ws.Range("B5").Value = 20
$ws.Range("B6").Value = 20
$ws.Range("B7").Value = 20
$ws.Range("B8").Value = 20
$ws.Range("B9").Value = 20
$ws.Range("B10").Value = 20
$ws.Range("B11").Value = 20
$ws.Range("B12").Value = 20
$ws.Range("B13").Value = 20
$ws.Range("B14").Value = 50
$ws.Range("B15").Value = 20

# Update selection to reflect new "next empty row" position
$ws.Range("A16:XFD16").Select()
